$d = $word.ActiveDocument

# =====================================================================
# Pass 1: append the new 'Knärot' section paragraphs with plain text.
# =====================================================================
$anchor = $d.Paragraphs.Last.Range

# Paragraph 1
$anchor.InsertParagraphAfter()
$p0 = $d.Paragraphs.Last
$p0.Style = "Heading1"
$p0.Range.Text = "Knärot – ekologi samt krav på livsmiljön"
$p0Start = $p0.Range.Start
$anchor = $p0.Range

# Paragraph 2
$anchor.InsertParagraphAfter()
$p1 = $d.Paragraphs.Last
$p1.Style = "Normal"
$p1.Range.Text = "Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021)."
$p1Start = $p1.Range.Start
$anchor = $p1.Range

# Paragraph 3
$anchor.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Style = "Normal"
$p2.Range.Text = "Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”"
$p2Start = $p2.Range.Start
$anchor = $p2.Range

# Paragraph 4
$anchor.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Style = "Normal"
$p3.Range.Text = "Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”"
$p3Start = $p3.Range.Start
$anchor = $p3.Range

# Paragraph 5
$anchor.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last
$p4.Style = "Normal"
$p4.Range.Text = "En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022)."
$p4Start = $p4.Range.Start
$anchor = $p4.Range

# Paragraph 6
$anchor.InsertParagraphAfter()
$p5 = $d.Paragraphs.Last
$p5.Style = "Normal"
$p5.Range.Text = "Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022)."
$p5Start = $p5.Range.Start
$anchor = $p5.Range

# Paragraph 7
$anchor.InsertParagraphAfter()
$p6 = $d.Paragraphs.Last
$p6.Style = "Heading2"
$p6.Range.Text = "Referenser - knärot"
$p6Start = $p6.Range.Start
$anchor = $p6.Range

# Paragraph 8
$anchor.InsertParagraphAfter()
$p7 = $d.Paragraphs.Last
$p7.Style = "Normal"
$p7.Range.Text = "de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025"
$p7Start = $p7.Range.Start
$anchor = $p7.Range

# Paragraph 9
$anchor.InsertParagraphAfter()
$p8 = $d.Paragraphs.Last
$p8.Style = "Normal"
$p8.Range.Text = "Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 "
$p8Start = $p8.Range.Start
$anchor = $p8.Range

# Paragraph 10
$anchor.InsertParagraphAfter()
$p9 = $d.Paragraphs.Last
$p9.Style = "Normal"
$p9.Range.Text = "Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853"
$p9Start = $p9.Range.Start
$anchor = $p9.Range

# Paragraph 11
$anchor.InsertParagraphAfter()
$p10 = $d.Paragraphs.Last
$p10.Style = "Normal"
$p10.Range.Text = "Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62."
$p10Start = $p10.Range.Start
$anchor = $p10.Range

# Paragraph 12
$anchor.InsertParagraphAfter()
$p11 = $d.Paragraphs.Last
$p11.Style = "Normal"
$p11.Range.Text = "Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/"
$p11Start = $p11.Range.Start
$anchor = $p11.Range

# Paragraph 13
$anchor.InsertParagraphAfter()
$p12 = $d.Paragraphs.Last
$p12.Style = "Normal"
$p12.Range.Text = "SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala "
$p12Start = $p12.Range.Start
$anchor = $p12.Range

# =====================================================================
# Pass 2: apply italics to the relevant runs via absolute offsets.
# =====================================================================
# Paragraph 3 italics
$d.Range($p2Start + 34, $p2Start + 116).Font.Italic = 1
$d.Range($p2Start + 278, $p2Start + 483).Font.Italic = 1
$d.Range($p2Start + 490, $p2Start + 608).Font.Italic = 1

# Paragraph 4 italics
$d.Range($p3Start + 205, $p3Start + 1070).Font.Italic = 1

# Paragraph 8 italics
$d.Range($p7Start + 33, $p7Start + 113).Font.Italic = 1

# Paragraph 9 italics
$d.Range($p8Start + 62, $p8Start + 176).Font.Italic = 1

# Paragraph 10 italics
$d.Range($p9Start + 117, $p9Start + 207).Font.Italic = 1

# Paragraph 11 italics
$d.Range($p10Start + 54, $p10Start + 121).Font.Italic = 1

# Paragraph 12 italics
$d.Range($p11Start + 22, $p11Start + 57).Font.Italic = 1

# Paragraph 13 italics
$d.Range($p12Start + 25, $p12Start + 61).Font.Italic = 1

# =====================================================================
# Update the date stamp in the first-page header.
# =====================================================================
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(2)
$hdr.Range.Find.Execute("2023-09-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15", 2)

